$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (Polarity), shifting Polarity -> C and Review -> D
$ws.Columns("B:B").Insert()

# New column header
$ws.Range("B1").Value = "Unnamed: 0.1"

# New column values mirror column A's index values
$ws.Range("B2").Value = 0
$ws.Range("B3").Value = 1
$ws.Range("B4").Value = 2
$ws.Range("B5").Value = 3
$ws.Range("B6").Value = 4
